$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 437652
$ws.Range("B2").Value = 99121

# I2 keeps its original "text that looks numeric" storage (t="inlineStr" in
# the source) -- force Text formatting for the write, then drop back to the
# default style so no stray style index is left behind on the cell.
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "500"
$ws.Range("I2").Style = "Normal"

$ws.Range("P2").Value = "SSV Ängarna lokal 2, Dls"
$ws.Range("S2").Value = 10
$ws.Range("X2").ClearContents()

$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2007-05-31"
$ws.Range("Y2").Style = "Normal"

$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2007-05-31"
$ws.Range("AA2").Style = "Normal"

$ws.Range("AC2").Value = "10 m²"
$ws.Range("AI2").ClearContents()
$ws.Range("AW2").Value = "Lars Sjögren"
$ws.Range("AX2").Value = "Gunnar Flygh"
$ws.Range("AY2").Value = ""
